$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(30).Insert()

$ws.Cells.Item(30, 1).Value = "Social"
$ws.Cells.Item(30, 2).Value = "x"
$ws.Cells.Item(30, 8).Value = "Age45_54"
$ws.Cells.Item(30, 9).Value = "Total population between 45 and 54 years of age."
$ws.Cells.Item(30, 10).Value = 1980
$ws.Cells.Item(30, 11).Value = "IPUMS NHGIS"
$ws.Cells.Item(30, 12).Value = "Integrated Public Use Microdata Series National Historic Geographic Information System"
$ws.Cells.Item(30, 14).Value = "Integer"
$ws.Cells.Item(30, 17).Value = "The 1980 Census data does not allow for disaggregation into Age45_49 and Age50_54 variables, so this variable is used instead. "

$ws.Range("Q31").Select()
